$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$priceData = @{
    2 = @(18.94002225752671, 21.17923736572266, 21.52373157655332, 18.60270518053397, 3953196953)
    3 = @(20.83508725843124, 21.19668197631836, 21.62336385377512, 20.74830529858798, 3953196953)
    4 = @(20.70414632262139, 18.86037635803223, 21.28715587387261, 16.78340396294699, 3953196953)
    5 = @(21.20905738678914, 20.01893997192383, 21.20905738678914, 18.96840454920285, 3953196953)
    6 = @(17.36894796947403, 19.38274002075196, 19.92320568834162, 16.62858357700609, 3953196953)
    7 = @(20.53127499995558, 21.70427703857422, 21.72669160764143, 19.28355915616021, 3953196953)
    8 = @(23.14600198327657, 23.70391845703125, 23.89994339722005, 22.66347850049294, 3953196953)
    9 = @(23.45132074835026, 22.66834259033203, 24.24189935506609, 22.46309567756309, 3953196953)
    10 = @(23.65268541298835, 26.20579528808594, 26.31313277529562, 23.3230043167145, 3953196953)
    11 = @(26.37995125663196, 24.38463401794434, 26.75890521531853, 23.48751473078885, 3953196953)
    12 = @(24.65928530747341, 25.14325904846192, 25.3462177957519, 23.69914255680671, 3953196953)
    13 = @(26.9982183096425, 29.36813926696777, 29.94290411039217, 26.51005939865164, 3953196953)
    14 = @(32.59350694472726, 35.5204963684082, 36.40097367353828, 29.62685398226136, 3953196953)
    15 = @(35.11039457278382, 34.15087509155273, 37.07740752665912, 33.99895236836392, 3953196953)
    16 = @(33.83685647208515, 38.49456024169922, 38.54291122467185, 33.20024902928949, 3953196953)
    17 = @(37.05111522153125, 38.83592987060547, 39.29835871264406, 35.76118441621196, 3953196953)
    18 = @(38.72786106005871, 42.32513046264648, 42.48046591932168, 37.85307109911647, 3953196953)
    19 = @(46.07327796804423, 42.80701065063477, 46.75615127297581, 42.00895779143495, 3953196953)
    20 = @(46.23156700035268, 38.75536727905273, 47.1008920351993, 38.08474348504571, 3953196953)
    21 = @(39.08227709308084, 37.79761123657227, 40.85912436445526, 37.0718547900178, 3953196953)
    22 = @(38.99143897313606, 33.55448532104492, 42.25192757545226, 32.22675810235826, 3953196953)
    23 = @(35.39108580237228, 40.56571578979492, 40.56571578979492, 34.34767693222673, 3953196953)
    24 = @(40.30387063138773, 36.09735870361328, 41.4238987722378, 35.43902185088453, 3953196953)
    25 = @(31.22721489297705, 37.12060928344727, 37.14649428933779, 30.77852280523962, 3953196953)
    26 = @(39.07440315503297, 39.03091049194336, 42.91921496480492, 38.84823931600049, 3953196953)
    27 = @(44.86522085848997, 46.3458366394043, 47.43220217380526, 43.89274627487723, 3953196953)
    28 = @(49.05092186293772, 52.06808471679688, 53.1708482781522, 48.15989038612424, 3953196953)
    29 = @(49.62208612965672, 48.70727157592773, 52.0734386684144, 45.36775190301782, 3953196953)
    30 = @(49.5392642754138, 49.82515335083008, 51.54048780332775, 47.08240269050936, 3953196953)
    31 = @(44.3824310029588, 40.52347946166992, 46.46932249189246, 36.89840564304387, 3953196953)
    32 = @(40.89728591734684, 40.58869552612305, 45.38091951944369, 40.4253239388699, 3953196953)
    33 = @(41.75352049609955, 45.54596328735352, 45.64672924949458, 39.74737069025502, 3953196953)
    34 = @(44.68178799535382, 44.70948791503906, 47.77507367299224, 42.87198711359912, 3953196953)
    35 = @(43.91035005218153, 46.20819854736328, 46.97105110170853, 42.38465204114739, 3953196953)
    36 = @(49.2769524239324, 53.75752639770508, 53.92625126979812, 48.96762587558856, 3953196953)
    37 = @(49.13245131687603, 45.67694091796875, 50.51087815543385, 43.61874036108011, 3953196953)
    38 = @(47.79259432127118, 46.02285003662109, 48.15415587921645, 45.52808254406165, 3953196953)
    39 = @(44.79690720105008, 44.60505676269531, 48.06794455061136, 43.98154375285301, 3953196953)
    40 = @(47.0382292550883, 48.89576721191406, 49.36015077847237, 43.05226754085176, 3953196953)
    41 = @(53.36385449345381, 57.72147750854492, 58.36488515821725, 53.15913100403318, 3953196953)
    42 = @(59.07668407709886, 62.92417907714844, 65.26997143648015, 58.44852254504295, 3953196953)
    43 = @(57.08988588876331, 62.28706359863281, 64.96469498942184, 56.99108177666165, 3953196953)
    44 = @(67.10480483359422, 68.67531585693359, 72.11455549417215, 65.35536628936002, 3953196953)
}

foreach ($row in $priceData.Keys) {
    $vals = $priceData[$row]
    $ws.Cells.Item($row, 4).Value = $vals[0]
    $ws.Cells.Item($row, 5).Value = $vals[1]
    $ws.Cells.Item($row, 6).Value = $vals[2]
    $ws.Cells.Item($row, 7).Value = $vals[3]
    $ws.Cells.Item($row, 8).Value = $vals[4]
    $ws.Cells.Item($row, 9).Value = "CSCO"
}
